$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Game_Record")

# Append a new game record (row 15): 1st=SiderFace, 2nd=Doanage, 3rd=DrSystomatix, 4th=SimpleJack
$ws.Cells.Item(15, 1).Formula = "=ROW()-1"
$ws.Cells.Item(15, 2).Value = 45951
$ws.Cells.Item(15, 3).Value = "SiderFace"
$ws.Cells.Item(15, 4).Value = "Doanage"
$ws.Cells.Item(15, 5).Value = "DrSystomatix"
$ws.Cells.Item(15, 6).Value = "SimpleJack"

# Move the sheet selection to reflect where the user clicked next
[void]$ws.Range("G15").Select()
